$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Widen column J (LastLoginDate) slightly to fit the refreshed timestamps
$ws.Columns.Item(10).ColumnWidth = 14.448721

# Update LastLoginDate values to reflect the latest login timestamps
$ws.Range("J2").Value = 45966.4738045602
$ws.Range("J3").Value = 45966.5204188657
$ws.Range("J4").Value = 45966.5186435185
